# Updated cryptos list with GitHub Actions - applies the latest scraped
# coinranking.com price / 1h-volume snapshot to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) hold text that frequently *looks* numeric
# ("314.01", "1.950", "116.00", ...). Excel's COM layer auto-coerces a
# plain .Value assignment like that into a real number (dropping
# significant trailing zeros / punctuation), so force the cell to Text
# first, assign, then restore the default "Normal" style so the cell
# ends up identical in shape to its neighbours (no stray style index).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Row 25 / 26: Monero and Toncoin swap ranking order ------------------
# Row 25 was Monero, becomes Toncoin; Row 26 was Toncoin, becomes Monero.
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D25") "1.950"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D26") "157.03"
$ws.Range("E26").Value = "  +2.57%  "

# --- Row 45 / 46: Decentraland and EnergySwap swap ranking order ---------
# Row 45 was Decentraland, becomes EnergySwap; Row 46 was EnergySwap, becomes Decentraland.
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "10.77"
$ws.Range("E45").Value = "  +1.57%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D46") "0.4899"
$ws.Range("E46").Value = "  +3.21%  "

# --- Price / volume refresh for all other rows ----------------------------
Set-TextValue $ws.Range("D2") "27.903.67"
$ws.Range("E2").Value = "  +2.78%  "

Set-TextValue $ws.Range("D3") "1.873.52"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").Value = "  -0.46%  "

Set-TextValue $ws.Range("D5") "314.01"
$ws.Range("E5").Value = "  +1.23%  "

Set-TextValue $ws.Range("D6") "1.014"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("E7").Value = "  +1.16%  "

Set-TextValue $ws.Range("D8") "0.3827"
$ws.Range("E8").Value = "  +3.55%  "

Set-TextValue $ws.Range("D9") "0.07383"
$ws.Range("E9").Value = "  +1.46%  "

Set-TextValue $ws.Range("D10") "0.9406"
$ws.Range("E10").Value = "  +0.84%  "

Set-TextValue $ws.Range("D11") "21.05"
$ws.Range("E11").Value = "  +5.60%  "

Set-TextValue $ws.Range("D12") "0.07820"
$ws.Range("E12").Value = "  +0.25%  "

Set-TextValue $ws.Range("D13") "1.875.05"
$ws.Range("E13").Value = "  +1.06%  "

Set-TextValue $ws.Range("D14") "5.508"
$ws.Range("E14").Value = "  +2.23%  "

Set-TextValue $ws.Range("D15") "6.613"
$ws.Range("E15").Value = "  +1.74%  "

Set-TextValue $ws.Range("D16") "91.32"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("E17").Value = "  -0.50%  "

Set-TextValue $ws.Range("D18") "0.000008881"
$ws.Range("E18").Value = "  +2.11%  "

Set-TextValue $ws.Range("D19") "1.012"
$ws.Range("E19").Value = "  -0.40%  "

Set-TextValue $ws.Range("D20") "28.023.93"
$ws.Range("E20").Value = "  +3.10%  "

Set-TextValue $ws.Range("D21") "14.88"
$ws.Range("E21").Value = "  +1.76%  "

Set-TextValue $ws.Range("D22") "5.130"
$ws.Range("E22").Value = "  +1.09%  "

Set-TextValue $ws.Range("D23") "2.123.28"
$ws.Range("E23").Value = "  +2.18%  "

Set-TextValue $ws.Range("D24") "10.85"
$ws.Range("E24").Value = "  +1.79%  "

Set-TextValue $ws.Range("D27") "18.58"
$ws.Range("E27").Value = "  +1.00%  "

Set-TextValue $ws.Range("D28") "2.067"
$ws.Range("E28").Value = "  +4.13%  "

Set-TextValue $ws.Range("D29") "116.00"
$ws.Range("E29").Value = "  +0.89%  "

Set-TextValue $ws.Range("D30") "4.990"
$ws.Range("E30").Value = "  +1.25%  "

Set-TextValue $ws.Range("D31") "0.08910"
$ws.Range("E31").Value = "  +0.42%  "

Set-TextValue $ws.Range("D32") "3.333"
$ws.Range("E32").Value = "  +0.31%  "

Set-TextValue $ws.Range("D33") "1.224"
$ws.Range("E33").Value = "  +3.72%  "

Set-TextValue $ws.Range("D34") "0.7652"
$ws.Range("E34").Value = "  +3.75%  "

Set-TextValue $ws.Range("D35") "4.658"
$ws.Range("E35").Value = "  +2.74%  "

Set-TextValue $ws.Range("D36") "2.721"
$ws.Range("E36").Value = "  +1.37%  "

Set-TextValue $ws.Range("D37") "1.135"
$ws.Range("E37").Value = "  +1.57%  "

Set-TextValue $ws.Range("D38") "0.02055"
$ws.Range("E38").Value = "  +3.37%  "

Set-TextValue $ws.Range("D39") "0.5651"
$ws.Range("E39").Value = "  +6.89%  "

Set-TextValue $ws.Range("D40") "0.05373"
$ws.Range("E40").Value = "  +1.96%  "

Set-TextValue $ws.Range("D41") "2.997"
$ws.Range("E41").Value = "  +0.50%  "

Set-TextValue $ws.Range("D42") "7.056"
$ws.Range("E42").Value = "  +0.14%  "

Set-TextValue $ws.Range("D43") "8.592"
$ws.Range("E43").Value = "  +3.39%  "

Set-TextValue $ws.Range("D44") "0.1537"
$ws.Range("E44").Value = "  +0.70%  "

Set-TextValue $ws.Range("D47") "105.40"
$ws.Range("E47").Value = "  +3.12%  "

$ws.Range("E48").Value = "  -0.43%  "

Set-TextValue $ws.Range("D49") "1.676"
$ws.Range("E49").Value = "  +3.21%  "

Set-TextValue $ws.Range("D50") "67.78"
$ws.Range("E50").Value = "  +2.61%  "

Set-TextValue $ws.Range("D51") "0.06113"
$ws.Range("E51").Value = "  +0.91%  "
